# Data update from DGS's 2021/09/22 report: append a new row (84) to the
# risk-matrix time series with the date 2021/09/22 and its associated
# incidence / R(t) figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 84

# Column A holds the date as plain text (it is displayed via a
# "yyyy/mm/dd" number format, but the underlying cell type in this sheet
# is always text/shared-string, matching every other row above it).
# Temporarily switch the cell to a Text format so Excel's COM layer does
# not silently convert the typed string into a date serial number, then
# restore the original date display format once the text value is set.
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("A$newRow").Value = "2021/09/22"
$ws.Range("A$newRow").NumberFormat = "yyyy/mm/dd"

# Numeric columns: Incidência Portugal, Incidência Continente, Rt Portugal, Rt Continente
$ws.Range("B$newRow").Value = 137.4
$ws.Range("C$newRow").Value = 140.1
$ws.Range("D$newRow").Value = 0.82
$ws.Range("E$newRow").Value = 0.81

# Move the active selection to the next empty row, as happens naturally
# after entering a new row of data.
$ws.Range("A85").Select() | Out-Null
